$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like
# "18.00" or "8.10" keep their trailing zeros and formats like
# "2.351.19" are not misinterpreted as numbers/dates.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '60.174.03'
$ws.Range('E2').Value = '  +3.51%  '
$ws.Range('D3').Value = '2.351.19'
$ws.Range('E3').Value = '  +2.90%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '549.64'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('D6').Value = '134.27'
$ws.Range('E6').Value = '  +2.67%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '0.589'
$ws.Range('E8').Value = '  +3.36%  '
$ws.Range('D9').Value = '2.345.61'
$ws.Range('E9').Value = '  +2.79%  '
$ws.Range('D10').Value = '0.102'
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('D11').Value = '5.56'
$ws.Range('E11').Value = '  +2.08%  '
$ws.Range('D12').Value = '0.151'
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('D13').Value = '0.337'
$ws.Range('E13').Value = '  +2.64%  '
$ws.Range('D14').Value = '24.22'
$ws.Range('E14').Value = '  +3.27%  '
$ws.Range('D15').Value = '2.768.73'
$ws.Range('E15').Value = '  +2.87%  '
$ws.Range('D16').Value = '59.972.71'
$ws.Range('E16').Value = '  +3.25%  '
$ws.Range('D17').Value = '0.0000134'
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('D18').Value = '2.343.38'
$ws.Range('E18').Value = '  +2.71%  '
$ws.Range('D19').Value = '10.73'
$ws.Range('E19').Value = '  +1.87%  '
$ws.Range('D20').Value = '4.23'
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').Value = '317.25'
$ws.Range('E21').Value = '  +1.79%  '
$ws.Range('D22').Value = '6.73'
$ws.Range('E22').Value = '  +5.32%  '
$ws.Range('D23').Value = '1.01'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').Value = '63.32'
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('D25').Value = '0.175'
$ws.Range('E25').Value = '  +4.60%  '
$ws.Range('D26').Value = '0.996'
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('D27').Value = '8.10'
$ws.Range('E27').Value = '  +2.16%  '
$ws.Range('D28').Value = '1.34'
$ws.Range('E28').Value = '  +5.48%  '
$ws.Range('D29').Value = '1.75'
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0752'
$ws.Range('E30').Value = '  +5.49%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '171.25'
$ws.Range('E31').Value = '  +0.63%  '
$ws.Range('B32').Value = 'SuiNetwork'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D32').Value = '1.17'
$ws.Range('E32').Value = '  +8.10%  '
$ws.Range('D33').Value = '5.94'
$ws.Range('E33').Value = '  +4.31%  '
$ws.Range('D34').Value = '1.44'
$ws.Range('E34').Value = '  +17.73%  '
$ws.Range('D35').Value = '0.388'
$ws.Range('E35').Value = '  +2.47%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').Value = '18.00'
$ws.Range('E36').Value = '  +2.05%  '
$ws.Range('B37').Value = 'USDe'
$ws.Range('C37').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D39').Value = '4.15'
$ws.Range('E39').Value = '  +6.73%  '
$ws.Range('D40').Value = '320.23'
$ws.Range('E40').Value = '  +12.58%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '1.55'
$ws.Range('E41').Value = '  +4.82%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').Value = '38.35'
$ws.Range('E42').Value = '  +0.99%  '
$ws.Range('D43').Value = '145.03'
$ws.Range('E43').Value = '  +4.55%  '
$ws.Range('D44').Value = '3.48'
$ws.Range('E44').Value = '  +2.73%  '
$ws.Range('D45').Value = '0.0959'
$ws.Range('E45').Value = '  +1.38%  '
$ws.Range('D46').Value = '0.0501'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').Value = '0.566'
$ws.Range('E47').Value = '  +2.95%  '
$ws.Range('D48').Value = '18.81'
$ws.Range('E48').Value = '  +4.13%  '
$ws.Range('D49').Value = '0.0213'
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = '1.55'
$ws.Range('E51').Value = '  +4.87%  '
